$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3a6421e250>),
                ('model',
                 BaggingClassifier(estimator=DecisionTreeClassifier(criterion='entropy',
                                                                    max_depth=5,
                                                                    max_features='sqrt',
                                                                    min_samples_split=6,
                                                                    random_state=42),
                                   n_estimators=50, random_state=42))])"
$ws.Range("B2").Value = 0.6571428571428571
$ws.Range("C2").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f3a5c78c340>, 'scaler': None, 'model__n_estimators': 50, 'model__estimator__min_samples_split': 6, 'model__estimator__min_samples_leaf': 1, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 5, 'model__estimator__criterion': 'entropy', 'model__estimator__class_weight': None}"
$ws.Range("D2").Value = 0.1818181818181818
$ws.Range("E2").Value = "[1 1 0 0 1 0 0 0 0 1 0 1]"
$ws.Range("F2").Value = "[0 0 1 0 0 1 1 1 1 1 0 0]"
$ws.Range("G2").Value = 77
$ws.Range("H2").Value = 0.8387582881253767
$ws.Range("I2").Value = 0.02611026980688669
$ws.Range("J2").Value = 0.5780590717299577
$ws.Range("K2").Value = 0.06679224815257635

$ws.Range("A3").Value = "Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3a64669400>),
                ('model',
                 BaggingClassifier(estimator=DecisionTreeClassifier(criterion='entropy',
                                                                    max_depth=4,
                                                                    max_features='log2',
                                                                    min_samples_leaf=4,
                                                                    random_state=42),
                                   n_estimators=50, random_state=42))])"
$ws.Range("B3").Value = 0.6666666666666666
$ws.Range("C3").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f45a9c6ad00>, 'scaler': None, 'model__n_estimators': 50, 'model__estimator__min_samples_split': 2, 'model__estimator__min_samples_leaf': 4, 'model__estimator__max_features': 'log2', 'model__estimator__max_depth': 4, 'model__estimator__criterion': 'entropy', 'model__estimator__class_weight': None}"
$ws.Range("D3").Value = 0.7777777777777778
$ws.Range("E3").Value = "[1 1 0 1 0 0 1 0 1 1 1 0]"
$ws.Range("F3").Value = "[1 1 1 1 1 1 1 1 1 1 1 0]"
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 0.837092731829574
$ws.Range("I3").Value = 0.02737973134630577
$ws.Range("J3").Value = 0.5185463659147869
$ws.Range("K3").Value = 0.08552962744456431

$ws.Range("A4").Value = "Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3a646692e0>),
                ('model',
                 BaggingClassifier(estimator=DecisionTreeClassifier(max_depth=2,
                                                                    max_features='sqrt',
                                                                    min_samples_leaf=5,
                                                                    random_state=42),
                                   random_state=42))])"
$ws.Range("B4").Value = 0.6095238095238095
$ws.Range("C4").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f45a9c6aa60>, 'scaler': None, 'model__n_estimators': 10, 'model__estimator__min_samples_split': 2, 'model__estimator__min_samples_leaf': 5, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 2, 'model__estimator__criterion': 'gini', 'model__estimator__class_weight': None}"
$ws.Range("D4").Value = 0.8421052631578948
$ws.Range("E4").Value = "[1 0 1 1 1 1 0 1 0 1 0 1]"
$ws.Range("F4").Value = "[1 1 1 1 1 1 0 1 1 1 1 1]"
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.8395833333333333
$ws.Range("I4").Value = 0.02733684835345362
$ws.Range("J4").Value = 0.5041005291005292
$ws.Range("K4").Value = 0.08684466418134143

$ws.Range("A5").Value = "Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f45a99f0fd0>),
                ('model',
                 BaggingClassifier(estimator=DecisionTreeClassifier(max_depth=2,
                                                                    max_features='sqrt',
                                                                    min_samples_split=3,
                                                                    random_state=42),
                                   random_state=42))])"
$ws.Range("B5").Value = 0.6476190476190475
$ws.Range("C5").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f3a6464e4c0>, 'scaler': None, 'model__n_estimators': 10, 'model__estimator__min_samples_split': 3, 'model__estimator__min_samples_leaf': 1, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 2, 'model__estimator__criterion': 'gini', 'model__estimator__class_weight': None}"
$ws.Range("D5").Value = 0.6666666666666666
$ws.Range("E5").Value = "[1 1 0 0 0 0 1 0 1 1 1 1]"
$ws.Range("F5").Value = "[1 0 0 0 1 1 1 1 1 0 1 1]"
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 0.8398001175778953
$ws.Range("I5").Value = 0.02455359392544976
$ws.Range("J5").Value = 0.5243974132863022
$ws.Range("K5").Value = 0.07578724510946509

$ws.Range("A6").Value = "Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3a64669a90>),
                ('model',
                 BaggingClassifier(estimator=DecisionTreeClassifier(max_depth=6,
                                                                    max_features='log2',
                                                                    min_samples_leaf=6,
                                                                    min_samples_split=3,
                                                                    random_state=42),
                                   n_estimators=5, random_state=42))])"
$ws.Range("B6").Value = 0.6190476190476191
$ws.Range("C6").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f3a6464e040>, 'scaler': None, 'model__n_estimators': 5, 'model__estimator__min_samples_split': 3, 'model__estimator__min_samples_leaf': 6, 'model__estimator__max_features': 'log2', 'model__estimator__max_depth': 6, 'model__estimator__criterion': 'gini', 'model__estimator__class_weight': None}"
$ws.Range("D6").Value = 0.7142857142857143
$ws.Range("E6").Value = "[1 1 1 1 0 0 0 0 1 1 0 0]"
$ws.Range("F6").Value = "[1 1 1 0 0 1 0 1 1 1 1 0]"
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 0.8331863609641388
$ws.Range("I6").Value = 0.02723023566380105
$ws.Range("J6").Value = 0.5334509112286889
$ws.Range("K6").Value = 0.08353027919085367

Write-Host "done"